$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserList")

# Update the password hash for cmoticska (row 1, column H)
$ws.Range("H1").Value = "5e884898da28047151d0e56f8dc6292773603d0d6aabbdd62a11ef721d1542d8"
$ws.Range("H1").NumberFormat = "0.00E+00"

# Select the H1 cell (matches diff's new selection)
$ws.Range("H1").Select()
